$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start from a clean sheet so the workbook's shared-string / style tables
# only contain what we actually need.
$ws.Cells.Clear()

# ---------------------------------------------------------------------
# Pass 1 : write every cell value. A leading apostrophe is used for the
# cells that must end up with quotePrefix="1" in the saved file (text
# that looks like a number / date / formula-ish string).
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "New Delhi"
$ws.Range("A1").Value = "From"
$ws.Range("B1").Value = "To"
$ws.Range("C1").Value = "'DepartureDate"
$ws.Range("D1").Value = "'FareOption"
$ws.Range("D2").Value = "SAVER"
$ws.Range("D3").Value = "'FLEXI PLUS"
$ws.Range("C3").Value = "'25 April 2024"
$ws.Range("H1").Value = "Traveller1.Gender"
$ws.Range("I1").Value = "Traveller1.Country"
$ws.Range("E1").Value = "NoOfTraveller"
$ws.Range("E2").Value = "'1"
$ws.Range("F1").Value = "Traveller1.FirstMiddleName"
$ws.Range("G1").Value = "Traveller1.LastName"
$ws.Range("F2").Value = "Anwar"
$ws.Range("G2").Value = "Khan"
$ws.Range("H2").Value = "MALE"
$ws.Range("I2").Value = "India"
$ws.Range("J1").Value = "MobileNo"
$ws.Range("K1").Value = "Email"
$ws.Range("L1").Value = "CountryCode"
$ws.Range("J2").Value = "'7358101855"
$ws.Range("K2").Value = "'abcd123@gmail.com"
$ws.Range("C2").Value = "'10 January 2025"
$ws.Range("A2").Value = "Pune"
$ws.Range("B2").Value = "New Delhi"
$ws.Range("L2").Value = "India"
$ws.Range("A3").Value = "New Delhi"
$ws.Range("B3").Value = "Pune"
$ws.Range("E3").Value = "'1"

# ---------------------------------------------------------------------
# Pass 2 : bold the header row. A1 first so the "bold only" style is
# created before the "bold + quotePrefix" style used by C1/D1.
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("C1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true
$ws.Range("G1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true
$ws.Range("I1").Font.Bold = $true
$ws.Range("J1").Font.Bold = $true
$ws.Range("K1").Font.Bold = $true
$ws.Range("L1").Font.Bold = $true

# ---------------------------------------------------------------------
# Pass 3 : date formatting for the two DepartureDate values.
# ---------------------------------------------------------------------
$ws.Range("C2").NumberFormat = "d-mmm-yy"
$ws.Range("C3").NumberFormat = "d-mmm-yy"

# ---------------------------------------------------------------------
# Column widths for the newly added columns D:L.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 12.65
$ws.Columns.Item(5).ColumnWidth = 12.65
$ws.Columns.Item(6).ColumnWidth = 29.55
$ws.Columns.Item(7).ColumnWidth = 29.55
$ws.Columns.Item(8).ColumnWidth = 19.05
$ws.Columns.Item(9).ColumnWidth = 16.55
$ws.Columns.Item(10).ColumnWidth = 19.7

# ---------------------------------------------------------------------
# Sheet view / print setup.
# ---------------------------------------------------------------------
$ws.Range("F9").Select()
$ws.PageSetup.Orientation = 1
